$wb = $excel.ActiveWorkbook

# --- Sheets involved ---
$wsBeneficiaries = $wb.Worksheets.Item("Beneficiaries_allTiers")
$wsDisability    = $wb.Worksheets.Item("Disability_allTiers")

# --- Beneficiaries sheet: content + layout updates ---
# B5 held a stray/unrecognized placeholder string; replace with the correct
# variable name for this sheet.
$wsBeneficiaries.Range("B5").Value = "nbeneficiaries"

# Give column B an explicit width (previously default width).
$wsBeneficiaries.Columns.Item(2).ColumnWidth = 95 / 7

# Move the stored selection on the Beneficiaries sheet.
$wsBeneficiaries.Range("D24").Select()

# --- Disability sheet: rename + becomes the active/selected tab ---
$wsDisability.Name = "Disb_allTiers"
$wsDisability.Activate()
